# Update res_bus vm_pu results for Case_3_24 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value2 = 1.02
$ws.Cells.Item(2, 3).Value2 = 1.056680474601468
$ws.Cells.Item(2, 4).Value2 = 1.063156687288193
$ws.Cells.Item(2, 5).Value2 = 1.063321178223701
$ws.Cells.Item(2, 6).Value2 = 1.074901414098865
$ws.Cells.Item(2, 9).Value2 = 1.050713376927185
$ws.Cells.Item(2, 10).Value2 = 1.061680846436217
$ws.Cells.Item(2, 11).Value2 = 1.065876088249257
$ws.Cells.Item(2, 12).Value2 = 1.066040134039477
$ws.Cells.Item(2, 13).Value2 = 1.077589401630622
$ws.Cells.Item(2, 14).Value2 = 1.024283042918393
$ws.Cells.Item(3, 2).Value2 = 1.02
$ws.Cells.Item(3, 3).Value2 = 1.057671467213246
$ws.Cells.Item(3, 4).Value2 = 1.063950906507674
$ws.Cells.Item(3, 5).Value2 = 1.064189142821644
$ws.Cells.Item(3, 6).Value2 = 1.075816918658263
$ws.Cells.Item(3, 9).Value2 = 1.050991991129252
$ws.Cells.Item(3, 10).Value2 = 1.062323809802936
$ws.Cells.Item(3, 11).Value2 = 1.066485315557876
$ws.Cells.Item(3, 12).Value2 = 1.06672295364499
$ws.Cells.Item(3, 13).Value2 = 1.078321877910922
$ws.Cells.Item(3, 14).Value2 = 1.024501637001314
$ws.Cells.Item(4, 2).Value2 = 1.02
$ws.Cells.Item(4, 3).Value2 = 1.058313229880396
$ws.Cells.Item(4, 4).Value2 = 1.064465279099526
$ws.Cells.Item(4, 5).Value2 = 1.064751564188039
$ws.Cells.Item(4, 6).Value2 = 1.076410145136492
$ws.Cells.Item(4, 9).Value2 = 1.051171388272868
$ws.Cells.Item(4, 10).Value2 = 1.062739737547908
$ws.Cells.Item(4, 11).Value2 = 1.066879327129088
$ws.Cells.Item(4, 12).Value2 = 1.067164929440821
$ws.Cells.Item(4, 13).Value2 = 1.078796036208186
$ws.Cells.Item(4, 14).Value2 = 1.024642939986186
$ws.Cells.Item(5, 2).Value2 = 1.02
$ws.Cells.Item(5, 3).Value2 = 1.05858315153497
$ws.Cells.Item(5, 4).Value2 = 1.064681630122715
$ws.Cells.Item(5, 5).Value2 = 1.064988194043055
$ws.Cells.Item(5, 6).Value2 = 1.076659735838593
$ws.Cells.Item(5, 9).Value2 = 1.051246594474671
$ws.Cells.Item(5, 10).Value2 = 1.062914565729982
$ws.Cells.Item(5, 11).Value2 = 1.067044920957239
$ws.Cells.Item(5, 12).Value2 = 1.067350769978704
$ws.Cells.Item(5, 13).Value2 = 1.078995418462022
$ws.Cells.Item(5, 14).Value2 = 1.024702309367839
$ws.Cells.Item(6, 2).Value2 = 1.02
$ws.Cells.Item(6, 3).Value2 = 1.058628479851921
$ws.Cells.Item(6, 4).Value2 = 1.064717962754279
$ws.Cells.Item(6, 5).Value2 = 1.065027936219321
$ws.Cells.Item(6, 6).Value2 = 1.076701654805112
$ws.Cells.Item(6, 9).Value2 = 1.051259209459187
$ws.Cells.Item(6, 10).Value2 = 1.062943918492331
$ws.Cells.Item(6, 11).Value2 = 1.067072722010328
$ws.Cells.Item(6, 12).Value2 = 1.067381975382309
$ws.Cells.Item(6, 13).Value2 = 1.079028898291617
$ws.Cells.Item(6, 14).Value2 = 1.024712275719886
$ws.Cells.Item(7, 2).Value2 = 1.02
$ws.Cells.Item(7, 3).Value2 = 1.058316836098109
$ws.Cells.Item(7, 4).Value2 = 1.064468169566604
$ws.Cells.Item(7, 5).Value2 = 1.064754725310964
$ws.Cells.Item(7, 6).Value2 = 1.076413479402178
$ws.Cells.Item(7, 9).Value2 = 1.051172394016876
$ws.Cells.Item(7, 10).Value2 = 1.062742073721152
$ws.Cells.Item(7, 11).Value2 = 1.066881539994204
$ws.Cells.Item(7, 12).Value2 = 1.067167412519181
$ws.Cells.Item(7, 13).Value2 = 1.078798700183825
$ws.Cells.Item(7, 14).Value2 = 1.024643733418407
$ws.Cells.Item(8, 2).Value2 = 1.02
$ws.Cells.Item(8, 3).Value2 = 1.05701527613995
$ws.Cells.Item(8, 4).Value2 = 1.063425001463481
$ws.Cells.Item(8, 5).Value2 = 1.063614346536102
$ws.Cells.Item(8, 6).Value2 = 1.075210639707813
$ws.Cells.Item(8, 9).Value2 = 1.050807718809315
$ws.Cells.Item(8, 10).Value2 = 1.061898161679331
$ws.Cells.Item(8, 11).Value2 = 1.066082020158347
$ws.Cells.Item(8, 12).Value2 = 1.066270865318219
$ws.Cells.Item(8, 13).Value2 = 1.077836903935206
$ws.Cells.Item(8, 14).Value2 = 1.024356946976098
$ws.Cells.Item(9, 2).Value2 = 1.02
$ws.Cells.Item(9, 3).Value2 = 1.054725813850599
$ws.Cells.Item(9, 4).Value2 = 1.06159038178578
$ws.Cells.Item(9, 5).Value2 = 1.061610961793768
$ws.Cells.Item(9, 6).Value2 = 1.073097527022051
$ws.Cells.Item(9, 9).Value2 = 1.05015836027081
$ws.Cells.Item(9, 10).Value2 = 1.060410258531516
$ws.Cells.Item(9, 11).Value2 = 1.064671679640514
$ws.Cells.Item(9, 12).Value2 = 1.064692196005017
$ws.Cells.Item(9, 13).Value2 = 1.076143654524518
$ws.Cells.Item(9, 14).Value2 = 1.023850522573919
$ws.Cells.Item(10, 2).Value2 = 1.02
$ws.Cells.Item(10, 3).Value2 = 1.053202280309503
$ws.Cells.Item(10, 4).Value2 = 1.060369780196602
$ws.Cells.Item(10, 5).Value2 = 1.060279553781436
$ws.Cells.Item(10, 6).Value2 = 1.071693196305571
$ws.Cells.Item(10, 9).Value2 = 1.04972094382439
$ws.Cells.Item(10, 10).Value2 = 1.059417822916556
$ws.Cells.Item(10, 11).Value2 = 1.06373050982645
$ws.Cells.Item(10, 12).Value2 = 1.063640589717172
$ws.Cells.Item(10, 13).Value2 = 1.075015933938505
$ws.Cells.Item(10, 14).Value2 = 1.023512210246773
$ws.Cells.Item(11, 2).Value2 = 1.02
$ws.Cells.Item(11, 3).Value2 = 1.052543241622474
$ws.Cells.Item(11, 4).Value2 = 1.059841849406636
$ws.Cells.Item(11, 5).Value2 = 1.059704046831899
$ws.Cells.Item(11, 6).Value2 = 1.071086168066576
$ws.Cells.Item(11, 9).Value2 = 1.049530474044087
$ws.Cells.Item(11, 10).Value2 = 1.05898797954601
$ws.Cells.Item(11, 11).Value2 = 1.063322760696779
$ws.Cells.Item(11, 12).Value2 = 1.06318544401363
$ws.Cells.Item(11, 13).Value2 = 1.074527895670464
$ws.Cells.Item(11, 14).Value2 = 1.023365557023962
$ws.Cells.Item(12, 2).Value2 = 1.02
$ws.Cells.Item(12, 3).Value2 = 1.052298545425799
$ws.Cells.Item(12, 4).Value2 = 1.059645843422153
$ws.Cells.Item(12, 5).Value2 = 1.059490429617683
$ws.Cells.Item(12, 6).Value2 = 1.070860850715792
$ws.Cells.Item(12, 9).Value2 = 1.049459565469169
$ws.Cells.Item(12, 10).Value2 = 1.058828300557671
$ws.Cells.Item(12, 11).Value2 = 1.063171272888628
$ws.Cells.Item(12, 12).Value2 = 1.063016414648362
$ws.Cells.Item(12, 13).Value2 = 1.074346658464358
$ws.Cells.Item(12, 14).Value2 = 1.023311059563185
$ws.Cells.Item(13, 2).Value2 = 1.02
$ws.Cells.Item(13, 3).Value2 = 1.052351029048101
$ws.Cells.Item(13, 4).Value2 = 1.059687883233005
$ws.Cells.Item(13, 5).Value2 = 1.059536244351215
$ws.Cells.Item(13, 6).Value2 = 1.070909174793082
$ws.Cells.Item(13, 9).Value2 = 1.049474782814621
$ws.Cells.Item(13, 10).Value2 = 1.058862552954291
$ws.Cells.Item(13, 11).Value2 = 1.063203768954657
$ws.Cells.Item(13, 12).Value2 = 1.06305267055698
$ws.Cells.Item(13, 13).Value2 = 1.074385532545796
$ws.Cells.Item(13, 14).Value2 = 1.02332275053007
$ws.Cells.Item(14, 2).Value2 = 1.02
$ws.Cells.Item(14, 3).Value2 = 1.052523012911843
$ws.Cells.Item(14, 4).Value2 = 1.059825645627083
$ws.Cells.Item(14, 5).Value2 = 1.059686386064778
$ws.Cells.Item(14, 6).Value2 = 1.071067539993428
$ws.Cells.Item(14, 9).Value2 = 1.049524615972576
$ws.Cells.Item(14, 10).Value2 = 1.058974780753019
$ws.Cells.Item(14, 11).Value2 = 1.063310239307347
$ws.Cells.Item(14, 12).Value2 = 1.063171471334974
$ws.Cells.Item(14, 13).Value2 = 1.074512913684933
$ws.Cells.Item(14, 14).Value2 = 1.023361052732713
$ws.Cells.Item(15, 2).Value2 = 1.02
$ws.Cells.Item(15, 3).Value2 = 1.052628991135262
$ws.Cells.Item(15, 4).Value2 = 1.05991053768171
$ws.Cells.Item(15, 5).Value2 = 1.05977891346709
$ws.Cells.Item(15, 6).Value2 = 1.071165135255292
$ws.Cells.Item(15, 9).Value2 = 1.049555298689167
$ws.Cells.Item(15, 10).Value2 = 1.059043925905799
$ws.Cells.Item(15, 11).Value2 = 1.063375835041288
$ws.Cells.Item(15, 12).Value2 = 1.063244672685931
$ws.Cells.Item(15, 13).Value2 = 1.074591402997122
$ws.Cells.Item(15, 14).Value2 = 1.023384648825016
$ws.Cells.Item(16, 2).Value2 = 1.02
$ws.Cells.Item(16, 3).Value2 = 1.053246032535178
$ws.Cells.Item(16, 4).Value2 = 1.060404829921771
$ws.Cells.Item(16, 5).Value2 = 1.060317769487506
$ws.Cells.Item(16, 6).Value2 = 1.071733505133212
$ws.Cells.Item(16, 9).Value2 = 1.049733562266271
$ws.Cells.Item(16, 10).Value2 = 1.059446347912683
$ws.Cells.Item(16, 11).Value2 = 1.063757566309589
$ws.Cells.Item(16, 12).Value2 = 1.063670800675838
$ws.Cells.Item(16, 13).Value2 = 1.075048329254726
$ws.Cells.Item(16, 14).Value2 = 1.023521939762445
$ws.Cells.Item(17, 2).Value2 = 1.02
$ws.Cells.Item(17, 3).Value2 = 1.053633263902979
$ws.Cells.Item(17, 4).Value2 = 1.060715047275115
$ws.Cells.Item(17, 5).Value2 = 1.060656048629486
$ws.Cells.Item(17, 6).Value2 = 1.072090312225895
$ws.Cells.Item(17, 9).Value2 = 1.049845097346996
$ws.Cells.Item(17, 10).Value2 = 1.059698747155148
$ws.Cells.Item(17, 11).Value2 = 1.063996958938424
$ws.Cells.Item(17, 12).Value2 = 1.063938155531724
$ws.Cells.Item(17, 13).Value2 = 1.075335020632478
$ws.Cells.Item(17, 14).Value2 = 1.02360801570087
$ws.Cells.Item(18, 2).Value2 = 1.02
$ws.Cells.Item(18, 3).Value2 = 1.053859193172354
$ws.Cells.Item(18, 4).Value2 = 1.060896049402034
$ws.Cells.Item(18, 5).Value2 = 1.060853457678228
$ws.Cells.Item(18, 6).Value2 = 1.072298533591221
$ws.Cells.Item(18, 9).Value2 = 1.049910051028413
$ws.Cells.Item(18, 10).Value2 = 1.059845956481164
$ws.Cells.Item(18, 11).Value2 = 1.064136571670775
$ws.Cells.Item(18, 12).Value2 = 1.064094118981171
$ws.Cells.Item(18, 13).Value2 = 1.075502269054729
$ws.Cells.Item(18, 14).Value2 = 1.023658206707354
$ws.Cells.Item(19, 2).Value2 = 1.02
$ws.Cells.Item(19, 3).Value2 = 1.053936240007458
$ws.Cells.Item(19, 4).Value2 = 1.060957776181764
$ws.Cells.Item(19, 5).Value2 = 1.060920785401347
$ws.Cells.Item(19, 6).Value2 = 1.072369548928263
$ws.Cells.Item(19, 9).Value2 = 1.049932181096747
$ws.Cells.Item(19, 10).Value2 = 1.059896149185992
$ws.Cells.Item(19, 11).Value2 = 1.064184172426322
$ws.Cells.Item(19, 12).Value2 = 1.06414730183096
$ws.Cells.Item(19, 13).Value2 = 1.075559300855544
$ws.Cells.Item(19, 14).Value2 = 1.023675317880558
$ws.Cells.Item(20, 2).Value2 = 1.02
$ws.Cells.Item(20, 3).Value2 = 1.053591711041076
$ws.Cells.Item(20, 4).Value2 = 1.060681757923069
$ws.Cells.Item(20, 5).Value2 = 1.060619744483657
$ws.Cells.Item(20, 6).Value2 = 1.07205201966016
$ws.Cells.Item(20, 9).Value2 = 1.049833141322624
$ws.Cells.Item(20, 10).Value2 = 1.059671668224622
$ws.Cells.Item(20, 11).Value2 = 1.063971276534237
$ws.Cells.Item(20, 12).Value2 = 1.063909468822189
$ws.Cells.Item(20, 13).Value2 = 1.075304258644143
$ws.Cells.Item(20, 14).Value2 = 1.023598782181461
$ws.Cells.Item(21, 2).Value2 = 1.02
$ws.Cells.Item(21, 3).Value2 = 1.052472365180831
$ws.Cells.Item(21, 4).Value2 = 1.059785075513411
$ws.Cells.Item(21, 5).Value2 = 1.059642168875036
$ws.Cells.Item(21, 6).Value2 = 1.071020900962973
$ws.Cells.Item(21, 9).Value2 = 1.04950994575204
$ws.Cells.Item(21, 10).Value2 = 1.058941732899053
$ws.Cells.Item(21, 11).Value2 = 1.063278887308114
$ws.Cells.Item(21, 12).Value2 = 1.06313648657548
$ws.Cells.Item(21, 13).Value2 = 1.0744754019483
$ws.Cells.Item(21, 14).Value2 = 1.023349774346081
$ws.Cells.Item(22, 2).Value2 = 1.02
$ws.Cells.Item(22, 3).Value2 = 1.051769167624243
$ws.Cells.Item(22, 4).Value2 = 1.059221822890486
$ws.Cells.Item(22, 5).Value2 = 1.059028406537458
$ws.Cells.Item(22, 6).Value2 = 1.07037352185652
$ws.Cells.Item(22, 9).Value2 = 1.049305816605784
$ws.Cells.Item(22, 10).Value2 = 1.058482701307661
$ws.Cells.Item(22, 11).Value2 = 1.062843372179898
$ws.Cells.Item(22, 12).Value2 = 1.062650667655449
$ws.Cells.Item(22, 13).Value2 = 1.073954509996742
$ws.Cells.Item(22, 14).Value2 = 1.023193075006595
$ws.Cells.Item(23, 2).Value2 = 1.02
$ws.Cells.Item(23, 3).Value2 = 1.052141890761509
$ws.Cells.Item(23, 4).Value2 = 1.059520363385884
$ws.Cells.Item(23, 5).Value2 = 1.059353689907698
$ws.Cells.Item(23, 6).Value2 = 1.070716621554401
$ws.Cells.Item(23, 9).Value2 = 1.049414116684595
$ws.Cells.Item(23, 10).Value2 = 1.058726051106946
$ws.Cells.Item(23, 11).Value2 = 1.063074263996557
$ws.Cells.Item(23, 12).Value2 = 1.062908191536663
$ws.Cells.Item(23, 13).Value2 = 1.074230621278933
$ws.Cells.Item(23, 14).Value2 = 1.023276157267817
$ws.Cells.Item(24, 2).Value2 = 1.02
$ws.Cells.Item(24, 3).Value2 = 1.053610486793749
$ws.Cells.Item(24, 4).Value2 = 1.060696799770498
$ws.Cells.Item(24, 5).Value2 = 1.060636148466261
$ws.Cells.Item(24, 6).Value2 = 1.07206932210783
$ws.Cells.Item(24, 9).Value2 = 1.049838544053201
$ws.Cells.Item(24, 10).Value2 = 1.059683904061412
$ws.Cells.Item(24, 11).Value2 = 1.063982881371969
$ws.Cells.Item(24, 12).Value2 = 1.06392243105006
$ws.Cells.Item(24, 13).Value2 = 1.075318158581699
$ws.Cells.Item(24, 14).Value2 = 1.023602954459765
$ws.Cells.Item(25, 2).Value2 = 1.02
$ws.Cells.Item(25, 3).Value2 = 1.055317209353449
$ws.Cells.Item(25, 4).Value2 = 1.062064243574082
$ws.Cells.Item(25, 5).Value2 = 1.062128153186672
$ws.Cells.Item(25, 6).Value2 = 1.073643045548366
$ws.Cells.Item(25, 9).Value2 = 1.050327032306785
$ws.Cells.Item(25, 10).Value2 = 1.060795009188912
$ws.Cells.Item(25, 11).Value2 = 1.06503645669858
$ws.Cells.Item(25, 12).Value2 = 1.065100176153199
$ws.Cells.Item(25, 13).Value2 = 1.0765812087573
$ws.Cells.Item(25, 14).Value2 = 1.023981569819025
